# Applies the commit: "Descobrindo comandos sem o Google e criando novos apelidos"
# - Turns the trailing empty list paragraph (ilvl 2) into an ilvl-1 heading-ish
#   bullet with text, then appends 8 further bullet paragraphs under it.

function Set-StdFont($range) {
    $range.Font.Name = "Times New Roman"
    $range.Font.NameBi = "Times New Roman"
    $range.Font.Size = 12
    $range.Font.SizeBi = 12
}

# Fill $paragraph (already created, empty, with the right list level) with the
# given ordered list of run-texts. Uses extra scratch paragraphs + merge for
# paragraphs that need more than two runs, because the COM shim here only
# keeps complex-script run formatting intact for the first two edits made
# directly inside one paragraph.
function Fill-Paragraph($doc, $paraIndex, $runTexts) {
    $n = $runTexts.Count

    # First (and, if present, second) run: edit the paragraph directly.
    $p = $doc.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $r.InsertAfter($runTexts[0])
    Set-StdFont $r

    if ($n -ge 2) {
        $p = $doc.Paragraphs.Item($paraIndex)
        $r = $p.Range
        $r.Collapse(0)
        $r.InsertAfter($runTexts[1])
        Set-StdFont $r
    }

    if ($n -le 2) {
        return
    }

    # Remaining runs: each gets its own scratch paragraph appended right
    # after the paragraph being built, then we merge everything back down
    # into one paragraph (deleting the paragraph marks) so the final result
    # is a single paragraph with N separate runs, in order.
    $lastParaIndex = $paraIndex
    for ($i = 2; $i -lt $n; $i++) {
        $prevPara = $doc.Paragraphs.Item($lastParaIndex)
        $prevPara.Range.InsertParagraphAfter()
        $lastParaIndex = $lastParaIndex + 1
        $scratch = $doc.Paragraphs.Item($lastParaIndex)
        $rs = $scratch.Range
        $rs.InsertAfter($runTexts[$i])
        Set-StdFont $rs
    }

    $extraParas = $n - 2
    for ($i = 0; $i -lt $extraParas; $i++) {
        $target = $doc.Paragraphs.Item($paraIndex)
        $rt = $target.Range
        $mark = $doc.Range($rt.End - 1, $rt.End)
        $mark.Delete()
    }
}

$d = $word.ActiveDocument

# The document ends with a trailing empty "PargrafodaLista" paragraph at
# ilvl=2 (ListLevelNumber 3). That is the paragraph the diff turns into the
# new sub-heading line.
$headingIndex = $d.Paragraphs.Count
$headingPara = $d.Paragraphs.Item($headingIndex)
$headingPara.Range.ListFormat.ListLevelNumber = 2

$newParagraphs = @(
    @{ ilvl = 2; runs = @(
        "Get-command: Lista todos os comandos existentes no PS."
    ) },
    @{ ilvl = 3; runs = @(
        " Para procurar um específico podemos ver os argumentos que podemos passar para ele e temos -Name, -Verbose, -Verb dentre outros. Como queremos achar o rename, podemos utilizar o -Verb, uma vez renomear é um verbo."
    ) },
    @{ ilvl = 3; runs = @(
        " Entretanto, colocamos o rename entre `u2018*`u2019, para que ele sirva como coringa e traga tudo o que achar que possua esse verbo."
    ) },
    @{ ilvl = 2; runs = @(
        "Isso só é possível pois a nova nomenclatura de comandos no PS é composta por verb-noun, ou seja, um verbo e um noum."
    ) },
    @{ ilvl = 2; runs = @(
        "Para obter ajuda sobre algum comando ou qualquer coisa no PS, nós utilizamos um comando que pede exatamente isso: Get-Help -Name comando",
        " e ele nos devolve uma lista e texto com a ajuda do comando pedido",
        "."
    ) },
    @{ ilvl = 2; runs = @(
        "Para voltar o tab caso tenha encontrado o que procura, utilizamos shift+tab."
    ) },
    @{ ilvl = 2; runs = @(
        "-WhatIf: é um switch argument",
        ", ou seja, não precisa de valor nenhum após ter sido colocado. Esse argumento não executará o comando que colocamos, mas sim nos dirá o que aconteceria se o executarmos para evitar que utilizemos comandos que não conhecemos direito e façamos alguma coisa de errado."
    ) },
    @{ ilvl = 2; runs = @(
        "Para setar um novo alias utilizamos o comando new-alias -Name `u201capelido_do_comando`u201d comando-em-si",
        ": ",
        "New-Alias -name `"ren`" Rename-Item"
    ) }
)

# First fill the heading paragraph's own text (3 runs: " ", title, ":").
Fill-Paragraph $d $headingIndex @(
    " ",
    "Descobrindo comandos sem o Google e criando novos apelidos",
    ":"
)

$prevIndex = $headingIndex
foreach ($entry in $newParagraphs) {
    $prevPara = $d.Paragraphs.Item($prevIndex)
    $prevPara.Range.InsertParagraphAfter()
    $newIndex = $prevIndex + 1
    $newPara = $d.Paragraphs.Item($newIndex)
    $newPara.Range.ListFormat.ListLevelNumber = $entry.ilvl + 1

    Fill-Paragraph $d $newIndex $entry.runs

    $prevIndex = $newIndex
}

Write-Output "done"
